# Daily attendance processing - 2025-12-05 23:26:29
# Reorders the "Recorded By" value in column G from
#   "dnasr281@gmail.com, System"  ->  "System, dnasr281@gmail.com"
# for a specific set of rows on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 6, 10, 12, 13, 14, 15, 18, 19, 20, 21, 22, 24, 26, 29, 32, 36, 38, 39, 40, 41, 44, 45, 46, 47, 48, 50, 52, 55, 58, 62, 64, 65, 66, 67, 70, 71, 72, 73, 74, 76, 78, 83, 84, 85, 86, 90, 92, 99, 101, 109, 110, 111, 112, 116, 118, 125, 127, 135, 136, 137, 138, 142, 144, 151, 153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
}
